# Refresh the cryptos table: Price (D) and Volume(1h) (E) columns updated
# with the latest scrape values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Price cells whose new text would otherwise be auto-coerced to a number --
# -- by COM ("218.27" etc. look numeric): force Text storage, matching the --
# -- original inline-string cells, then drop back to the default (General) --
# -- style so no extra formatting is introduced.                           --
$textPrices = [ordered]@{
    "D5" = "218.27"
    "D6" = "0.5190"
    "D8" = "0.2574"
    "D9" = "0.06402"
    "D10" = "19.91"
    "D11" = "0.07784"
    "D14" = "4.289"
    "D15" = "0.5528"
    "D17" = "64.34"
    "D20" = "210.85"
    "D21" = "4.381"
    "D23" = "5.905"
    "D25" = "143.75"
    "D26" = "1.763"
    "D27" = "0.1163"
    "D28" = "6.955"
    "D30" = "0.05270"
    "D32" = "3.360"
    "D33" = "3.217"
    "D34" = "1.572"
    "D35" = "2.765"
    "D37" = "0.9260"
    "D38" = "0.5709"
    "D40" = "0.01592"
    "D42" = "0.8401"
    "D43" = "5.670"
    "D44" = "99.87"
    "D47" = "0.4511"
    "D48" = "55.99"
    "D49" = "1.010"
    "D50" = "7.876"
    "D51" = "0.05090"
}
foreach ($addr in $textPrices.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textPrices[$addr]
    $cell.Style = "Normal"
}

# -- Price cells whose new text already fails numeric parsing (multiple
# -- separators), so COM leaves them as plain text automatically. --
$plainPrices = [ordered]@{
    "D2" = "26.199.81"
    "D3" = "1.658.90"
    "D12" = "1.664.71"
    "D13" = "1.887.50"
    "D16" = "0.0₅8037"
    "D18" = "26.230.53"
    "D39" = "1.160.59"
    "D45" = "1.797.15"
}
foreach ($addr in $plainPrices.Keys) {
    $ws.Range($addr).Value = $plainPrices[$addr]
}

# -- Volume(1h) percentage cells: always text (padded with spaces), safe to --
# -- assign directly.                                                       --
$volumes = [ordered]@{
    "E2" = "  -4.31%  "
    "E3" = "  -2.89%  "
    "E4" = "  +0.38%  "
    "E5" = "  -2.67%  "
    "E6" = "  -2.71%  "
    "E7" = "  +0.39%  "
    "E9" = "  -3.09%  "
    "E10" = "  -4.93%  "
    "E12" = "  -2.37%  "
    "E13" = "  -2.90%  "
    "E14" = "  -5.76%  "
    "E15" = "  -4.15%  "
    "E16" = "  -1.74%  "
    "E17" = "  -5.00%  "
    "E18" = "  -4.09%  "
    "E19" = "  +0.50%  "
    "E20" = "  -2.74%  "
    "E21" = "  -6.16%  "
    "E22" = "  -3.44%  "
    "E23" = "  -0.98%  "
    "E24" = "  +0.39%  "
    "E25" = "  +1.11%  "
    "E26" = "  +1.20%  "
    "E27" = "  -4.10%  "
    "E28" = "  -4.19%  "
    "E29" = "  -3.17%  "
    "E30" = "  -2.37%  "
    "E31" = "  -2.56%  "
    "E32" = "  -3.75%  "
    "E33" = "  -6.06%  "
    "E34" = "  -4.17%  "
    "E35" = "  -3.75%  "
    "E36" = "  -2.14%  "
    "E37" = "  -2.41%  "
    "E38" = "  -2.56%  "
    "E39" = "  +11.10%  "
    "E40" = "  -2.52%  "
    "E41" = "  +0.40%  "
    "E42" = "  +0.17%  "
    "E43" = "  -3.07%  "
    "E44" = "  -1.03%  "
    "E45" = "  -2.95%  "
    "E46" = "  -3.59%  "
    "E47" = "  -0.05%  "
    "E48" = "  -3.35%  "
    "E49" = "  +0.42%  "
    "E50" = "  -2.29%  "
    "E51" = "  -2.80%  "
}
foreach ($addr in $volumes.Keys) {
    $ws.Range($addr).Value = $volumes[$addr]
}
